{"js": "// Correciones casos de uso extendido\n// Split the sentence in the \"1.2\" table row so it describes that a link is\n// emailed to the user before redirecting them to the password-change tab.\n\nconst body = context.document.body;\n\nconst originalText =\n  \"Si el usuario agrega bien sus credenciales, lo dirigir\u00e1 a otra pesta\u00f1a donde podr\u00e1 hacer el cambio de contrase\u00f1a.\";\nconst newText =\n  \"Si el usuario agrega bien sus credenciales, le enviara un link al correo electr\u00f3nico que lo dirigir\u00e1 a otra pesta\u00f1a donde podr\u00e1 hacer el cambio de contrase\u00f1a.\";\n\nconst results = body.search(originalText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence not found in document body.\");\n}\n\n// Replace the whole sentence in place so the existing run formatting\n// (Arial, sz 22) carries over to the newly inserted text automatically.\nresults.items[0].insertText(newText, \"Replace\");\n\nawait context.sync();\n", "ps1": "# Correciones casos de uso extendido\n# Split the sentence in the \"1.2\" table row so it describes that a link is\n# emailed to the user before redirecting them to the password-change tab.\n\n$d = $word.ActiveDocument\n\n$originalText = \"Si el usuario agrega bien sus credenciales, lo dirigir\u00e1 a otra pesta\u00f1a donde podr\u00e1 hacer el cambio de contrase\u00f1a.\"\n$newText = \"Si el usuario agrega bien sus credenciales, le enviara un link al correo electr\u00f3nico que lo dirigir\u00e1 a otra pesta\u00f1a donde podr\u00e1 hacer el cambio de contrase\u00f1a.\"\n\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = $originalText\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif ($found) {\n    # $rng now spans exactly the matched sentence; overwrite it with the\n    # corrected wording while keeping the existing run formatting (Arial, 22).\n    $rng.Text = $newText\n} else {\n    Write-Output \"Target sentence not found in document content.\"\n}\n"}
